$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("clients")
$ws2 = $wb.Worksheets.Item("drivers")
$ws4 = $wb.Worksheets.Item("keys")

# --- clients: fix address typo (1254 -> 1250 Portage Ave) ---
$ws1.Range("D3").Value = "1250 Portage Ave, Winnipeg, MB R3G 0T6"

# --- move the last two driver rows (Hertha / Henrietta) from "drivers" to "clients" ---
$ws2.Range("A5:F6").Copy()
$ws1.Range("A7").PasteSpecial(-4104)
$excel.CutCopyMode = $false
$ws2.Range("A5:F6").Clear()
$ws2.Range("A5:F6").Select()

# --- keys sheet: remove the leaked API_KEY row, add blank styled row 7 ---
$ws4.Range("A3:B3").Clear()
$ws4.Range("B7").Font.Name = "Calibri"
$ws4.Range("C7").Font.Name = "Menlo"
$ws4.Range("B7:C7").Select()

# --- restore selection/active sheet to "clients" ---
$ws1.Range("B11").Select()
